# Fix: last packer plan + add skus
# Adds a new SKU row (row 8) to the brynza params sheet, matching the
# new product "Брынза классическая "Из Лавки", 45%, 0,2 кг, т/ф".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row data (row 8) ---
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 'Брынза классическая "Из Лавки", 45%, 0,2 кг, т/ф'
$ws.Range("C8").Value = 45
$ws.Range("D8").Value = "Брынза"
$ws.Range("E8").Value = "Брынза"
$ws.Range("F8").Value = "Из Лавки"
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0.2
$ws.Range("I8").Value = 6
$ws.Range("J8").Value = 350
$ws.Range("K8").Value = "00-00012173"

# --- Formatting to match the other data rows ---
# B column (product name) wraps its text just like a long SKU title.
$ws.Range("B8").WrapText = $true

# K column (code) uses the same "Calibri 11 black" text font used by the
# other code/line cells in the sheet (same look as K2/F3/etc.).
$ws.Range("K2").Copy()
$ws.Range("K8").PasteSpecial(-4122)
$ws.Range("K8").Value = "00-00012173"

# Row 8 ends up a little taller than the other rows because of the
# wrapped product name.
$ws.Rows.Item(8).RowHeight = 14.9

# Leave the selection where the editor last left it.
$ws.Range("K21").Select()
